{"js": "// Merge-field renames in the correspondence template:\n//   ${Date}                        -> ${currentDate}\n//   ${Employee Name}  (both spots) -> ${personFirstName}\n//   ${Case Number}                 -> ${caseNumber}\n//   ${Case Title}                  -> ${caseTitle}\n\nconst body = context.document.body;\n\n// ${Date} -> ${currentDate}\nlet results = body.search(\"${Date}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"${currentDate}\", \"Replace\");\n}\nawait context.sync();\n\n// \"Employee Name}\" -> \"personFirstName}\" (covers both \"Dear ${Employee Name},\"\n// and \"...inform you, ${Employee Name},\")\nresults = body.search(\"Employee Name}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"personFirstName}\", \"Replace\");\n}\nawait context.sync();\n\n// \"Case Number}\" -> \"caseNumber}\"\nresults = body.search(\"Case Number}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"caseNumber}\", \"Replace\");\n}\nawait context.sync();\n\n// \"Case Title}\" -> \"caseTitle}\"\nresults = body.search(\"Case Title}\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"caseTitle}\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Merge-field renames in the correspondence template:\n#   ${Date}                        -> ${currentDate}\n#   ${Employee Name}  (both spots) -> ${personFirstName}\n#   ${Case Number}                 -> ${caseNumber}\n#   ${Case Title}                  -> ${caseTitle}\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n\nReplace-All \"`${Date}\" \"`${currentDate}\"\nReplace-All \"Employee Name}\" \"personFirstName}\"\nReplace-All \"Case Number}\" \"caseNumber}\"\nReplace-All \"Case Title}\" \"caseTitle}\"\n"}
